$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.73%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8.31%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.083"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.83%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08051"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.43%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.923"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.12%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.939"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.17%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9301"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.43%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1460"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "13.21%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1935"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.28%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08982"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.59%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03503"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.47%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09800"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.87%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001399"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.25%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005906"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-5.63%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.738"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.61%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.186"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.50%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.85%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3462"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.62%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1309"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.89%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.810"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.46%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.37%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04359"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.15%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.12%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004278"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-11.35%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001300"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.10%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02069"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.63%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05059"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.70%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007447"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.04%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01009"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.49%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1351"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.23%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002140"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.32%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008939"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-9.61%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006185"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.15%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.10%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002786"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001598"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "27.69%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.10%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.10%"
